$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 263
$ws.Range("F5").Value = 317
$ws.Range("F6").Value = 456
$ws.Range("F8").Value = 2070
$ws.Range("F11").Value = 42
$ws.Range("F14").Value = 1344
$ws.Range("F15").Value = 57
$ws.Range("G18").Value = 78
$ws.Range("F19").Value = 14
$ws.Range("F22").Value = 150
$ws.Range("F23").Value = 7132
$ws.Range("F24").Value = 7764
$ws.Range("F36").Value = 1406
$ws.Range("F37").Value = 28
$ws.Range("F41").Value = 706
$ws.Range("F43").Value = 1359
$ws.Range("F45").Value = 235
$ws.Range("F49").Value = 152

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 292

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 180

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 180
$ws.Range("F7").Value = 317
$ws.Range("F9").Value = 456
$ws.Range("F10").Value = 2070
$ws.Range("F12").Value = 42
$ws.Range("F16").Value = 1344
$ws.Range("G17").Value = 78
$ws.Range("F18").Value = 14
$ws.Range("F21").Value = 150
$ws.Range("F23").Value = 7132
$ws.Range("F24").Value = 7764
$ws.Range("F32").Value = 1406
$ws.Range("F33").Value = 28
$ws.Range("F39").Value = 706
$ws.Range("F43").Value = 1359
$ws.Range("F45").Value = 235
$ws.Range("F47").Value = 152
$ws.Range("F49").Value = 292

